$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Ordem" values in column A (rows 2-15) with new order numbers ---
$newValues = @(
    685601397708,
    685601397707,
    685601397704,
    685601397703,
    685601397702,
    685601397701,
    685601397699,
    685601397698,
    685601397697,
    685601397696,
    685601397695,
    685601397694,
    685601397692,
    685601397691
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}

# --- Rows 16-28: the old order numbers were deleted, leaving the formatted cells empty ---
$ws.Range("A16:A28").ClearContents()

# --- New helper columns L & M (rows 8-21) get a date number format, no values ---
$ws.Range("L8").NumberFormat = "mm-dd-yy"
$ws.Range("L8").Copy()
$ws.Range("L8:M21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# --- Update the selection shown when the workbook is reopened ---
$null = $ws.Range("F6:V25").Select()
